$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.430.06'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.50%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.691.89'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.54%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.011'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.70%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.02'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.42%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5533'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +8.63%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.62%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2716'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.89%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06482'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.46%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.12'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.35%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07617'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.29%  '

# Row 12
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.689.85'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.79%  '

# Row 13
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.563'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.06%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5835'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.57%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008462'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.47%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.35'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.06%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.499.55'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.65%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.965'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.89%  '

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.56%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.98'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.90%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '190.75'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.58%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.258'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.22%  '

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.62%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '150.06'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.57%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1311'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +8.23%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.923'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +4.19%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.79'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.10%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06343'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -4.33%  '

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.32%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.594'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.30%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.593'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.35%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.678'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.07%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.048'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.17%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6252'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.89%  '

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.47%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.726'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.38%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.251'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.02%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.124.88'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.92%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01647'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.56%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8863'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.03%  '

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.75%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.74'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.78%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.842.98'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.61%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000110'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.23%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '57.59'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.37%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.267'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.00%  '

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.15%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05284'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.02%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4303'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.42%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.090'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.68%  '
